$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.241.05'
$ws.Range('E2').Value = '  +5.15%  '
$ws.Range('D3').Value = '2.709.18'
$ws.Range('E3').Value = '  +4.27%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '586.15'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.62%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '149.67'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.82%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +1.68%  '
$ws.Range('D9').Value = '2.736.33'
$ws.Range('E9').Value = '  +5.09%  '
$ws.Range('E10').Value = '  +2.99%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.113'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +7.61%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.388'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('E13').Value = '  +1.78%  '
$ws.Range('D14').Value = '3.192.11'
$ws.Range('E14').Value = '  +4.30%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '26.73'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +9.62%  '
$ws.Range('D16').Value = '63.109.55'
$ws.Range('E16').Value = '  +4.93%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.0000151'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +7.88%  '
$ws.Range('D18').Value = '2.724.39'
$ws.Range('E18').Value = '  +4.64%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.96'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +5.58%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.87'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +5.65%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '363.28'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +5.14%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '7.02'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.69%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.531'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.23%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '65.51'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.98%  '
$ws.Range('E26').Value = '  +4.08%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.66'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +8.33%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.995'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').Value = '0.0₃0865'
$ws.Range('E29').Value = '  +8.48%  '
$ws.Range('E30').Value = '  +6.26%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.08'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +11.01%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '169.99'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +1.83%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.21'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +24.04%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.17%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '20.55'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +5.88%  '
$ws.Range('E36').Value = '  +12.29%  '
$ws.Range('E37').Value = '  +8.26%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.82'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +11.06%  '
$ws.Range('E39').Value = '  +19.44%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '350.97'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +12.25%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.27'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +10.24%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '39.17'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.68%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.68'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +14.54%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '21.60'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +8.51%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0594'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +8.39%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '21.74'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +9.52%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '139.59'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0260'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +7.32%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.641'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +5.80%  '
$ws.Range('E50').Value = '  +1.81%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.144.12'
$ws.Range('E51').Value = '  +6.32%  '
